$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'67.940.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.12%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.524.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  +1.21%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'183.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.58%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.25%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +4.03%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "'7.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.32%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +1.92%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "'4.137.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.19%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'32.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.84%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.35%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "'67.895.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.07%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +0.34%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "'3.512.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +1.17%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'14.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.53%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'399.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.06%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'8.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.17%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "'73.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.08%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  +1.04%  "

# Row 24 - Dai
$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "

# Row 25 - LEO
$ws.Range("D25").Value = "'5.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.33%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +1.70%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "'10.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.63%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -1.14%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.33%  "

# Row 30 - Fetch.AI
$ws.Range("E30").Value = "  +1.78%  "

# Row 31 - NEARProtocol
$ws.Range("E31").Value = "  +0.19%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +0.81%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "'24.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.83%  "

# Row 34 - Aptos
$ws.Range("D34").Value = "'7.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.19%  "

# Row 35 - USDe
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "'1.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.52%  "

# Row 37 - Monero
$ws.Range("D37").Value = "'164.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.77%  "

# Row 38 - Mantle
$ws.Range("E38").Value = "  -1.85%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +1.86%  "

# Row 40 - RenderToken
$ws.Range("D40").Value = "'7.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.42%  "

# Row 41 - dogwifhat
$ws.Range("E41").Value = "  +6.37%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  +1.26%  "

# Row 43 - EnergySwap
$ws.Range("D43").Value = "'27.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.90%  "

# Row 44 - Maker
$ws.Range("D44").Value = "'2.882.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.69%  "

# Row 45 - Hedera (was InjectiveProtocol)
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0745"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "

# Row 46 - InjectiveProtocol (was Hedera)
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'27.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.80%  "

# Row 47 - OKB
$ws.Range("D47").Value = "'42.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "

# Row 48 - Bittensor
$ws.Range("D48").Value = "'350.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.81%  "

# Row 49 - VeChain
$ws.Range("D49").Value = "'0.0307"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.40%  "

# Row 50 - ONDO
$ws.Range("E50").Value = "  -0.65%  "

# Row 51 - Arweave
$ws.Range("D51").Value = "'33.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.80%  "
